$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (06-nov) before the old DK column ---
$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Columns("DK:DK").Insert()

$ws.Range("DK1").Value = "06-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 115).Value = "-"
}

# --- Sheet "Gaz": append the next day's closing price ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A143").NumberFormat = "@"
$wsGaz.Range("A143").Value = "2025-11-04"
$wsGaz.Range("A143").Style = "Normal"
$wsGaz.Range("B143").Value = 31.17

# --- Sheet "CO2": append the next day's closing price ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A143").NumberFormat = "@"
$wsCo2.Range("A143").Value = "2025-11-04"
$wsCo2.Range("A143").Style = "Normal"
$wsCo2.Range("B143").Value = 81.9
